# Apply "Add files via upload" edits to Self-Cross Pairs workbook
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------------
# Sheet2: fix genus name typo "Triofolium" -> "Trifolium"
# ---------------------------------------------------------------------------
$ws2.Range("A15").Value = "Trifolium"

# ---------------------------------------------------------------------------
# Sheet1: fix distribution note text (add missing clarifying detail)
# ---------------------------------------------------------------------------
$ws1.Range("S10").Value = "Europe (not N), Mediterranean Basin (Middle East, N Africa), W North America"

# ---------------------------------------------------------------------------
# Sheet2: add explanatory notes / citations in columns E (and F for row 3)
# ---------------------------------------------------------------------------
$ws2.Range("E3").Value = "Only 2 in dataset"
$ws2.Range("E3").Font.Bold = $true
$ws2.Range("F3").Value = "Likely to be a ploidy effect- exclude"

$ws2.Range("E4").Value = "Only 2 in dataset"
$ws2.Range("E4").Font.Bold = $true

$ws2.Range("E5").Value = "Biogeography and Phylogeny of Cardamine (Brassicaceae). Carlsen et al., 2009"

$ws2.Range("E6").Value = "Complex rearrangements are involved in Cephalanthera (Orchidaceae) chromosome evolution. Moscone et al., 2007."

$ws2.Range("E7").Value = "Phylogenetic relationships within Luzula DC. and Juncus L. (Juncaceae): A comparison of phylogenetic signals of trnL-trnF intergenic spacer, trnL intro and rbcL plastome sequence data. Drabkova et al., 2006."

$ws2.Range("E8").Value = "Chloroplast DNA characters, phylogeny and classification of Lathyrus (Fabaceae). Asmussen and Liston, 1998."

$ws2.Range("E9").Value = "Phylogenetic relationships within Luzula DC. and Juncus L. (Juncaceae): A comparison of phylogenetic signals of trnL-trnF intergenic spacer, trnL intro and rbcL plastome sequence data. Drabkova et al., 2006."

$ws2.Range("E10").Value = "Only 2 in dataset"
$ws2.Range("E10").Font.Bold = $true

$ws2.Range("E11").Value = "Phylogenetics of Papaver and Related Genera Based on DNA Sequences from ITS Nuclear Ribosomal RNA and Plastid trnL Intron and trnL-F Intergenic Spacers. Carolan et al., 2006."

$ws2.Range("E14").Value = "Tackling speciose genera: species composition and phylogenetic position of Senecio sect. Jacobaea (Asteraceae) based on plastid and nrDNA sequences. Pelser et al., 2002."

$ws2.Range("E15").Value = "Molecular phylogenetics of the clover genus (Trifolium- Leguminosae). Ellison et al., 2006"
$ws2.Range("E16").Value = "Molecular phylogenetics of the clover genus (Trifolium- Leguminosae). Ellison et al., 2006"
$ws2.Range("E17").Value = "Molecular phylogenetics of the clover genus (Trifolium- Leguminosae). Ellison et al., 2006"
$ws2.Range("E18").Value = "Molecular phylogenetics of the clover genus (Trifolium- Leguminosae). Ellison et al., 2006"

$ws2.Range("E19").Value = "Phylogeny of Veronica- a Combination of Molecular and Chemical Evidence. Taskova et al., 2004"

# ---------------------------------------------------------------------------
# Sheet1: unhide the supporting-data columns and autofit their widths
# ---------------------------------------------------------------------------
$ws1.Range("D1:D1,H1:Q1").EntireColumn.Hidden = $false
$ws1.Range("D1:D1,H1:Q1").EntireColumn.AutoFit() | Out-Null

# ---------------------------------------------------------------------------
# View state: make Sheet2 the active/front sheet with its own selection,
# and give Sheet1 a plain top-left view with a new selection.
# ---------------------------------------------------------------------------
$ws1.Range("C30").Select() | Out-Null
$ws2.Activate()
$ws2.Range("P2").Select() | Out-Null
